$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.235.02'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.578.98'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.47%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.50'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.05'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.582'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.74'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0998'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.52%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.29%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.033.34'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '58.233.98'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.63'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.585.48'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.43'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '335.02'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.04'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.85'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.419'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.156'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -6.06%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0727'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.24%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.01'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.89'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.89'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '36.83'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.825'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.73%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '285.68'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.589'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.65'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.72%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0535'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0946'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '18.43'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0226'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.907.40'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.81'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.38'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.04%  '
